$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4100.8613
$ws.Range("J17").Value = 4100.8613
$ws.Range("L17").Value = 12302.5839
$ws.Range("N17").Value = -12638.5839
$ws.Range("H28").Value = 265.16
$ws.Range("I28").Value = 246.61905
$ws.Range("J28").Value = 362.5
$ws.Range("K28").Value = 246.61905
$ws.Range("L28").Value = 362.5
$ws.Range("M28").Value = 238.38095
$ws.Range("N28").Value = -1332.5
$ws.Range("H33").Value = 198.78947
$ws.Range("I33").Value = 184.48215
$ws.Range("K33").Value = 184.48215
$ws.Range("M33").Value = 44.51785000000001
$ws.Range("H64").Value = 33806
$ws.Range("I64").Value = 85591
$ws.Range("J64").Value = 2735
$ws.Range("K64").Value = 85591
$ws.Range("L64").Value = 2735
$ws.Range("M64").Value = -85343
$ws.Range("N64").Value = -3231
$ws.Range("H67").Value = 33806
$ws.Range("I67").Value = 85591
$ws.Range("J67").Value = 2735
$ws.Range("K67").Value = 85591
$ws.Range("L67").Value = 2735
$ws.Range("M67").Value = -84733
$ws.Range("N67").Value = -4451
$ws.Range("H74").Value = 3011.9412
$ws.Range("I74").Value = 3169.4614
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 3169.4614
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -2233.4614
$ws.Range("N74").Value = -4372
$ws.Range("H75").Value = 41638.25
$ws.Range("J75").Value = 41638.25
$ws.Range("L75").Value = 41638.25
$ws.Range("N75").Value = -43510.25
$ws.Range("H77").Value = 3011.9412
$ws.Range("I77").Value = 3169.4614
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 15847.307
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -11167.307
$ws.Range("N77").Value = -21860
$ws.Range("H78").Value = 41638.25
$ws.Range("J78").Value = 41638.25
$ws.Range("L78").Value = 124914.75
$ws.Range("N78").Value = -134274.75
$ws.Range("H93").Value = 32014.154
$ws.Range("J93").Value = 32014.154
$ws.Range("L93").Value = 32014.154
$ws.Range("N93").Value = -37006.15399999999
$ws.Range("H95").Value = 31303.5
$ws.Range("J95").Value = 31303.5
$ws.Range("L95").Value = 31303.5
$ws.Range("N95").Value = -36795.5
$ws.Range("H106").Value = 39144.43
$ws.Range("I106").Value = 2125
$ws.Range("J106").Value = 88503.664
$ws.Range("K106").Value = 2125
$ws.Range("L106").Value = 88503.664
$ws.Range("M106").Value = -1494
$ws.Range("N106").Value = -89765.664

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 9264593
$ws.Range("I88").Value = 3600
$ws.Range("J88").Value = 15879587
$ws.Range("K88").Value = 3600
$ws.Range("L88").Value = 15879587
$ws.Range("M88").Value = -3194
$ws.Range("N88").Value = -15880399
$ws.Range("H91").Value = 9264593
$ws.Range("I91").Value = 3600
$ws.Range("J91").Value = 15879587
$ws.Range("K91").Value = 3600
$ws.Range("L91").Value = 15879587
$ws.Range("M91").Value = -2196
$ws.Range("N91").Value = -15882395

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2140
$ws.Range("I86").Value = 2140
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2140
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1017
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2140
$ws.Range("I89").Value = 2140
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10700
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -5084
$ws.Range("N89").ClearContents()
$ws.Range("H95").Value = 42616
$ws.Range("J95").Value = 42616
$ws.Range("L95").Value = 42616
$ws.Range("N95").Value = -48108
$ws.Range("H103").Value = 42664
$ws.Range("J103").Value = 42664
$ws.Range("L103").Value = 42664
$ws.Range("N103").Value = -45008

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 35773.375
$ws.Range("J28").Value = 35773.375
$ws.Range("L28").Value = 35773.375
$ws.Range("N28").Value = -36263.375
$ws.Range("H43").Value = 41049.668
$ws.Range("J43").Value = 41049.668
$ws.Range("L43").Value = 41049.668
$ws.Range("N43").Value = -41417.668
$ws.Range("H62").Value = 2999.8333
$ws.Range("I62").Value = 2999.6
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2999.6
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2375.6
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2999.8333
$ws.Range("I65").Value = 2999.6
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14998
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11878
$ws.Range("N65").Value = -21240
$ws.Range("H101").Value = 41049.668
$ws.Range("J101").Value = 41049.668
$ws.Range("L101").Value = 41049.668
$ws.Range("N101").Value = -47539.668
$ws.Range("H141").Value = 6997.6
$ws.Range("I141").Value = 6990
$ws.Range("J141").Value = 6999.5
$ws.Range("K141").Value = 6990
$ws.Range("L141").Value = 6999.5
$ws.Range("M141").Value = -1810
$ws.Range("N141").Value = -17359.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 7500
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -9496
$ws.Range("H78").Value = 2500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 22500
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -32484
$ws.Range("H113").Value = 5026.5654
$ws.Range("I113").Value = 6509.7646
$ws.Range("J113").Value = 824.1667
$ws.Range("K113").Value = 19529.2938
$ws.Range("L113").Value = 2472.5001
$ws.Range("M113").Value = -17359.2938
$ws.Range("N113").Value = -6812.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4654.788
$ws.Range("I70").Value = 4816.32
$ws.Range("J70").Value = 4150
$ws.Range("K70").Value = 4816.32
$ws.Range("L70").Value = 4150
$ws.Range("M70").Value = -4546.32
$ws.Range("N70").Value = -4690
$ws.Range("H73").Value = 4654.788
$ws.Range("I73").Value = 4816.32
$ws.Range("J73").Value = 4150
$ws.Range("K73").Value = 4816.32
$ws.Range("L73").Value = 4150
$ws.Range("M73").Value = -3880.32
$ws.Range("N73").Value = -6022
$ws.Range("H80").Value = 338080.12
$ws.Range("I80").Value = 560577.75
$ws.Range("J80").Value = 4333.6665
$ws.Range("K80").Value = 560577.75
$ws.Range("L80").Value = 4333.6665
$ws.Range("M80").Value = -559579.75
$ws.Range("N80").Value = -6329.6665
$ws.Range("H83").Value = 338080.12
$ws.Range("I83").Value = 560577.75
$ws.Range("J83").Value = 4333.6665
$ws.Range("K83").Value = 2802888.75
$ws.Range("L83").Value = 21668.3325
$ws.Range("M83").Value = -2797896.75
$ws.Range("N83").Value = -31652.3325

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000.5715
$ws.Range("I68").Value = 2533.6667
$ws.Range("J68").Value = 3350.75
$ws.Range("K68").Value = 2533.6667
$ws.Range("L68").Value = 3350.75
$ws.Range("M68").Value = -1784.6667
$ws.Range("N68").Value = -4848.75
$ws.Range("H71").Value = 3000.5715
$ws.Range("I71").Value = 2533.6667
$ws.Range("J71").Value = 3350.75
$ws.Range("K71").Value = 12668.3335
$ws.Range("L71").Value = 16753.75
$ws.Range("M71").Value = -8924.333500000001
$ws.Range("N71").Value = -24241.75
$ws.Range("H127").Value = 28180.5
$ws.Range("J127").Value = 50711
$ws.Range("L127").Value = 50711
$ws.Range("N127").Value = -60631

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5000751
$ws.Range("I14").Value = 1004
$ws.Range("J14").Value = 6667333.5
$ws.Range("K14").Value = 1004
$ws.Range("L14").Value = 6667333.5
$ws.Range("M14").Value = -836
$ws.Range("N14").Value = -6667669.5
$ws.Range("H68").Value = 29428.285
$ws.Range("J68").Value = 29428.285
$ws.Range("L68").Value = 29428.285
$ws.Range("N68").Value = -31050.285
$ws.Range("H69").Value = 23297.285
$ws.Range("J69").Value = 23297.285
$ws.Range("L69").Value = 23297.285
$ws.Range("N69").Value = -24795.285
$ws.Range("H71").Value = 29428.285
$ws.Range("J71").Value = 29428.285
$ws.Range("L71").Value = 88284.855
$ws.Range("N71").Value = -96396.855
$ws.Range("H72").Value = 23297.285
$ws.Range("J72").Value = 23297.285
$ws.Range("L72").Value = 69891.855
$ws.Range("N72").Value = -77379.855
$ws.Range("H80").Value = 37659.8
$ws.Range("J80").Value = 37659.8
$ws.Range("L80").Value = 37659.8
$ws.Range("N80").Value = -39655.8
$ws.Range("H83").Value = 37659.8
$ws.Range("J83").Value = 37659.8
$ws.Range("L83").Value = 112979.4
$ws.Range("N83").Value = -122963.4
$ws.Range("H101").Value = 36384
$ws.Range("I101").Value = 38000
$ws.Range("J101").Value = 35576
$ws.Range("K101").Value = 38000
$ws.Range("L101").Value = 35576
$ws.Range("M101").Value = -34755
$ws.Range("N101").Value = -42066
